$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Content change: the second paragraph ("This a second line for test") is
# removed entirely. The _GoBack bookmark that was anchored at the end of
# that paragraph survives and ends up attached to the (now only) remaining
# paragraph, right after "This is a test for git in word files".
#
# We reproduce that precisely the way a human editing in Word would: select
# the second paragraph's visible text (but not its paragraph mark, so the
# bookmark - which sits after the run but still inside that paragraph -
# stays put for a moment), delete it, and then delete the paragraph mark
# that separates the now-empty paragraph from the previous one. Deleting
# that mark merges the (bookmark-only) paragraph up into the first one,
# which is exactly how Word carries a trailing bookmark across a paragraph
# merge.
# ---------------------------------------------------------------------------

$needle = "This a second line for test"

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like ("*" + $needle + "*")) {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    # Delete the paragraph's text, leaving its paragraph mark (and the
    # bookmark that lives just before that mark) in place.
    $target = $d.Paragraphs($targetIndex)
    $textEnd = $target.Range.End - 1
    $textRange = $d.Range($target.Range.Start, $textEnd)
    $textRange.Delete()

    if ($targetIndex -gt 1) {
        # Re-fetch the preceding paragraph (indices/ranges shift after the
        # delete above) and remove its trailing paragraph mark, merging the
        # now-empty, bookmark-only paragraph into it.
        $prev = $d.Paragraphs($targetIndex - 1)
        $markRange = $d.Range($prev.Range.End - 1, $prev.Range.End)
        $markRange.Delete()
    }
}

# ---------------------------------------------------------------------------
# Style-exception cleanup: the document's latent style list also drops the
# "Normal Table", "Table Web 3" and "Table Theme" entries. Remove them via
# the Application.LatentStyles collection if this host exposes it (no-op,
# harmlessly, if it does not).
# ---------------------------------------------------------------------------

$staleLatentStyles = @("Normal Table", "Table Web 3", "Table Theme")
foreach ($name in $staleLatentStyles) {
    try {
        $word.LatentStyles($name).Delete()
    } catch {
    }
}
